$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SummaryReport")

# Row 2 updates
$ws.Range("A2").Value = "C099995"
$ws.Range("B2").Value = "appellant's opening brief"
$ws.Range("C2").Value = "Business Exception"
$ws.Range("D2").Value = "Document Processing Failure: Citations are less than 10 in the Doc. Case Number: C099995"

# Remove wrap text on D2 so it matches the default (General) style
$ws.Range("D2").ClearFormats()

# Row 3 updates
$ws.Range("A3").Value = "C100010"
$ws.Range("B3").Value = "appellant's opening brief"
$ws.Range("C3").Value = "Business Exception"
$ws.Range("D3").Value = "Document Processing Failure: Headings not found in the Document,Certificate of Compliance. Case Number: C100010"
$ws.Range("E3").Value = "Failed"
